$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I51").Value = 1750.5
$ws.Range("J51").Value = 4523.25
$ws.Range("K51").Value = 1750.5
$ws.Range("L51").Value = 4523.25
$ws.Range("M51").Value = -1266.5
$ws.Range("N51").Value = -5491.25

$ws.Range("H93").Value = 35733.5
$ws.Range("J93").Value = 35733.5
$ws.Range("L93").Value = 35733.5
$ws.Range("N93").Value = -40725.5

$ws.Range("H98").Value = 4805.8237
$ws.Range("I98").Value = 2849.9167
$ws.Range("J98").Value = 9500
$ws.Range("K98").Value = 2849.9167
$ws.Range("L98").Value = 9500
$ws.Range("M98").Value = -1351.9167
$ws.Range("N98").Value = -12496

$ws.Range("H113").Value = 19001.5
$ws.Range("I113").Value = 15000
$ws.Range("J113").Value = 20335.334
$ws.Range("K113").Value = 15000
$ws.Range("L113").Value = 20335.334
$ws.Range("M113").Value = -11746
$ws.Range("N113").Value = -26843.334

$ws.Range("H122").Value = 4805.8237
$ws.Range("I122").Value = 2849.9167
$ws.Range("J122").Value = 9500
$ws.Range("K122").Value = 8549.750100000001
$ws.Range("L122").Value = 28500
$ws.Range("M122").Value = -6099.750100000001
$ws.Range("N122").Value = -33400

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 549.1724
$ws.Range("I97").Value = 565.6539
$ws.Range("J97").Value = 406.33334
$ws.Range("K97").Value = 565.6539
$ws.Range("L97").Value = 406.33334
$ws.Range("M97").Value = -69.65390000000002
$ws.Range("N97").Value = -1398.33334

$ws.Range("H122").Value = 2183.1667
$ws.Range("I122").Value = 1212.7693
$ws.Range("J122").Value = 3330
$ws.Range("K122").Value = 3638.3079
$ws.Range("L122").Value = 9990
$ws.Range("M122").Value = -1188.3079
$ws.Range("N122").Value = -14890

$ws.Range("H139").Value = 42607.31
$ws.Range("J139").Value = 42607.31
$ws.Range("L139").Value = 42607.31
$ws.Range("N139").Value = -52887.31

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3500
$ws.Range("I99").Value = 2600
$ws.Range("J99").Value = 3757.1428
$ws.Range("K99").Value = 2600
$ws.Range("L99").Value = 3757.1428
$ws.Range("M99").Value = -1102
$ws.Range("N99").Value = -6753.1428

$ws.Range("H110").Value = 50000
$ws.Range("J110").Value = 50000
$ws.Range("L110").Value = 50000
$ws.Range("N110").Value = -58180

$ws.Range("H138").Value = 40340.668
$ws.Range("J138").Value = 40340.668
$ws.Range("L138").Value = 40340.668
$ws.Range("N138").Value = -50620.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1936.1642
$ws.Range("I58").Value = 1676.4237
$ws.Range("J58").Value = 3851.75
$ws.Range("K58").Value = 1676.4237
$ws.Range("L58").Value = 3851.75
$ws.Range("M58").Value = -1473.4237
$ws.Range("N58").Value = -4257.75

$ws.Range("H105").Value = 2582.8572
$ws.Range("I105").Value = 2445
$ws.Range("J105").Value = 2766.6667
$ws.Range("K105").Value = 2445
$ws.Range("L105").Value = 2766.6667
$ws.Range("M105").Value = -698
$ws.Range("N105").Value = -6260.6667

$ws.Range("H132").Value = 3188.889
$ws.Range("I132").Value = 1817.2941
$ws.Range("K132").Value = 5451.8823
$ws.Range("M132").Value = -2921.8823

$ws.Range("H136").Value = 1936.1642
$ws.Range("I136").Value = 1676.4237
$ws.Range("J136").Value = 3851.75
$ws.Range("K136").Value = 5029.2711
$ws.Range("L136").Value = 11555.25
$ws.Range("M136").Value = -2479.2711
$ws.Range("N136").Value = -16655.25

$ws.Range("H138").Value = 28794
$ws.Range("J138").Value = 28794
$ws.Range("L138").Value = 28794
$ws.Range("N138").Value = -39074

$ws.Range("H140").Value = 39446
$ws.Range("J140").Value = 39446
$ws.Range("L140").Value = 39446
$ws.Range("N140").Value = -49806

$ws.Range("H141").Value = 27035.715
$ws.Range("J141").Value = 27035.715
$ws.Range("L141").Value = 27035.715
$ws.Range("N141").Value = -37395.715

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 3833.3333
$ws.Range("J116").Value = 3833.3333
$ws.Range("L116").Value = 11499.9999
$ws.Range("N116").Value = -18383.9999

$ws.Range("H131").Value = 7693334
$ws.Range("J131").Value = 915.0492
$ws.Range("L131").Value = 2745.1476
$ws.Range("N131").Value = -12825.1476

$ws.Range("H132").Value = 1938.3125
$ws.Range("I132").Value = 828.5714
$ws.Range("J132").Value = 2801.4443
$ws.Range("K132").Value = 7457.1426
$ws.Range("L132").Value = 25212.9987
$ws.Range("M132").Value = -4927.1426
$ws.Range("N132").Value = -30272.9987

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 8547694
$ws.Range("I107").Value = 465
$ws.Range("J107").Value = 15873890
$ws.Range("K107").Value = 465
$ws.Range("L107").Value = 15873890
$ws.Range("M107").Value = 1455
$ws.Range("N107").Value = -15877730

$ws.Range("H113").Value = 1227.9333
$ws.Range("I113").Value = 1238.5
$ws.Range("J113").Value = 1215.8572
$ws.Range("K113").Value = 1238.5
$ws.Range("L113").Value = 1215.8572
$ws.Range("M113").Value = 931.5
$ws.Range("N113").Value = -5555.8572

$ws.Range("H122").Value = 3691.0527
$ws.Range("I122").Value = 1901.75
$ws.Range("J122").Value = 6758.4287
$ws.Range("K122").Value = 5705.25
$ws.Range("L122").Value = 20275.2861
$ws.Range("M122").Value = -3255.25
$ws.Range("N122").Value = -25175.2861

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 9587.625
$ws.Range("I40").Value = 9700.4
$ws.Range("J40").Value = 9536.362999999999
$ws.Range("K40").Value = 9700.4
$ws.Range("L40").Value = 9536.362999999999
$ws.Range("M40").Value = -9564.4
$ws.Range("N40").Value = -9808.362999999999

$ws.Range("H61").Value = 1600
$ws.Range("I61").Value = 1371.4286
$ws.Range("K61").Value = 1371.4286
$ws.Range("M61").Value = -1169.4286

$ws.Range("H113").Value = 1600
$ws.Range("I113").Value = 1371.4286
$ws.Range("K113").Value = 1371.4286
$ws.Range("M113").Value = 798.5714

$ws.Range("H122").Value = 8288.888999999999
$ws.Range("I122").Value = 3800
$ws.Range("J122").Value = 9571.429
$ws.Range("K122").Value = 11400
$ws.Range("L122").Value = 28714.287
$ws.Range("M122").Value = -8950
$ws.Range("N122").Value = -33614.287

$ws.Range("H139").Value = 43713.332
$ws.Range("J139").Value = 43713.332
$ws.Range("L139").Value = 43713.332
$ws.Range("N139").Value = -53993.332

$ws.Range("H140").Value = 69053.5
$ws.Range("J140").Value = 69053.5
$ws.Range("L140").Value = 69053.5
$ws.Range("N140").Value = -79413.5

$ws.Range("H141").Value = 40872.69
$ws.Range("J141").Value = 40872.69
$ws.Range("L141").Value = 40872.69
$ws.Range("N141").Value = -51232.69

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 38019.6
$ws.Range("I42").Value = 10000
$ws.Range("J42").Value = 45024.5
$ws.Range("K42").Value = 10000
$ws.Range("L42").Value = 45024.5
$ws.Range("M42").Value = -9622
$ws.Range("N42").Value = -45780.5

$ws.Range("H132").Value = 11496370
$ws.Range("I132").Value = 1162.6428
$ws.Range("J132").Value = 22225230
$ws.Range("K132").Value = 3487.9284
$ws.Range("L132").Value = 66675690
$ws.Range("M132").Value = -957.9284000000002
$ws.Range("N132").Value = -66680750

$ws.Range("H138").Value = 47043.75
$ws.Range("J138").Value = 47043.75
$ws.Range("L138").Value = 47043.75
$ws.Range("N138").Value = -57323.75

$ws.Range("H139").Value = 40052.5
$ws.Range("J139").Value = 40012.668
$ws.Range("L139").Value = 40012.668
$ws.Range("N139").Value = -50292.668

$ws.Range("H140").Value = 31038.166
$ws.Range("J140").Value = 31038.166
$ws.Range("L140").Value = 31038.166
$ws.Range("N140").Value = -41398.166

$ws.Range("H141").Value = 43276.332
$ws.Range("J141").Value = 43276.332
$ws.Range("L141").Value = 43276.332
$ws.Range("N141").Value = -53636.332

